# GSEA3/Notes.xlsx — "ajout de fichier de conversion Excel en xml"
#
# The sheet header row had the LastName/FirstName column labels swapped
# relative to the data that was actually entered under them (column B held
# the "FirstName" series Gsea3_FN*, column C held the "LastName" series
# Gsea3_LN*, yet the headers said the opposite). This pass:
#   - fixes the B1/C1 header labels so they match the data underneath,
#   - renumbers the CNE identifiers in column A,
#   - clears the (accidental/inconsistent) header font formatting on
#     A1:C11 so the whole block uses the plain default style, and
#   - leaves the selection on the A1:C11 block that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the swapped column headers (B1/C1) -------------------------------
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# --- renumber the CNE column (A2:A11) -------------------------------------
$ws.Range("A2").Value  = 17000021
$ws.Range("A3").Value  = 17000022
$ws.Range("A4").Value  = 17000023
$ws.Range("A5").Value  = 17000024
$ws.Range("A6").Value  = 17000025
$ws.Range("A7").Value  = 17000026
$ws.Range("A8").Value  = 17000027
$ws.Range("A9").Value  = 17000028
$ws.Range("A10").Value = 17000029
$ws.Range("A11").Value = 17000030

# --- drop the stray header formatting on the CNE/LastName/FirstName block -
$ws.Range("A1:C11").ClearFormats()

# --- leave the just-edited block selected, matching the saved view --------
$ws.Range("A1:C11").Select() | Out-Null
